# Ejercicios 52 al 59
# Mark exercises 52-59 (rows 53-60) as resolved ("si") in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 53; $row -le 60; $row++) {
    $ws.Cells.Item($row, 4).Value = "si"
}

# Reflect the cursor position left behind by the edit (row 61, column C).
$ws.Range("C61").Select() | Out-Null
